# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Source data refresh (gh-pages output regenerated at 456a3b4):
#   展览   (sheet 1): F5, F7, F8, F11, F14, F15, F18, F19, F28, F30, F31, F42, F49
#   全部类型 (sheet 4): F6, F8, F9, F12, F17, F18, F20, F21, F30, F32, F33, F42, F49

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

$exhibitionUpdates = @{
    5  = 212
    7  = 1087
    8  = 8452
    11 = 6983
    14 = 5156
    15 = 5156
    18 = 5682
    19 = 5682
    28 = 9481
    30 = 1748
    31 = 1197
    42 = 4944
    49 = 938
}

$allTypesUpdates = @{
    6  = 212
    8  = 1087
    9  = 8452
    12 = 6983
    17 = 5156
    18 = 5156
    20 = 5682
    21 = 5682
    30 = 9481
    32 = 1748
    33 = 1197
    42 = 4944
    49 = 938
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
